$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.92
$ws.Range("D6").Value = -7.897
$ws.Range("A14").Value = -21.676
$ws.Range("C15").Value = -12.932
$ws.Range("A16").Value = -22.085
$ws.Range("D18").Value = -8.675999999999998
$ws.Range("D19").Value = -7.85
$ws.Range("A21").Value = -20.216
$ws.Range("C21").Value = -12.739
$ws.Range("C22").Value = -12.5
$ws.Range("A23").Value = -19.991
$ws.Range("C24").Value = -12.061
$ws.Range("A25").Value = -21.839
$ws.Range("A26").Value = -21.737
$ws.Range("C27").Value = -12.563
$ws.Range("C28").Value = -12.596
$ws.Range("A29").Value = -21.279
$ws.Range("D35").Value = -8.270999999999999
$ws.Range("C36").Value = -13.045
$ws.Range("C39").Value = -12.82
$ws.Range("A40").Value = -19.934
$ws.Range("D44").Value = -7.448
$ws.Range("C45").Value = -12.678
$ws.Range("D47").Value = -7.494
$ws.Range("C48").Value = -11.416
$ws.Range("C49").Value = -12.884
$ws.Range("D50").Value = -8.683
$ws.Range("D51").Value = -8.093
$ws.Range("C52").Value = -11.382
$ws.Range("D52").Value = -7.535000000000001
$ws.Range("A53").Value = -21.977
$ws.Range("C53").Value = -12.81
$ws.Range("C54").Value = -12.348
$ws.Range("D55").Value = -8.451000000000001
$ws.Range("A57").Value = -21.802
$ws.Range("C57").Value = -11.657
$ws.Range("D57").Value = -8.013
$ws.Range("D58").Value = -8.407999999999998
$ws.Range("A59").Value = -22.266
$ws.Range("D64").Value = -7.657000000000001
$ws.Range("A65").Value = -21.484
$ws.Range("D66").Value = -7.755000000000001
$ws.Range("A69").Value = -21.484
$ws.Range("C70").Value = -11.664
$ws.Range("C71").Value = -11.081
$ws.Range("A79").Value = -21.117
$ws.Range("D80").Value = -7.983
$ws.Range("A83").Value = -22.259
$ws.Range("D83").Value = -8.218999999999999
$ws.Range("C86").Value = -13.468
$ws.Range("C87").Value = -13.462
$ws.Range("C89").Value = -13.376
$ws.Range("A91").Value = -20.744
$ws.Range("D92").Value = -7.717999999999999
$ws.Range("A93").Value = -21.508
$ws.Range("D94").Value = -7.782999999999999
$ws.Range("D96").Value = -7.712000000000001
$ws.Range("D97").Value = -8.023999999999999
$ws.Range("A100").Value = -22.323
$ws.Range("C101").Value = -12.375
$ws.Range("D101").Value = -7.659999999999999
$ws.Range("A103").Value = -22.072
